$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 307; this shifts rows 307:335 down to 308:336
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row 307 with the new weekly record
$ws.Cells.Item(307, 1).Value  = 8
$ws.Cells.Item(307, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(307, 3).Value  = "Coquimbo"
$ws.Cells.Item(307, 4).Value2 = 44578
$ws.Cells.Item(307, 5).Value  = 4
$ws.Cells.Item(307, 6).Value  = 100114001
$ws.Cells.Item(307, 7).Value  = "Papa"
$ws.Cells.Item(307, 8).Value  = "Asterix"
$ws.Cells.Item(307, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(307, 10).Value = 2000
$ws.Cells.Item(307, 11).Value = 10000
$ws.Cells.Item(307, 12).Value = 11000
$ws.Cells.Item(307, 13).Value = 10500
$ws.Cells.Item(307, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(307, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(307, 16).Value = 420
$ws.Cells.Item(307, 17).Value = 25
$ws.Cells.Item(307, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date style/number format used by the column (style index 2)
$ws.Cells.Item(307, 4).NumberFormat = $ws.Cells.Item(308, 4).NumberFormat
